$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reordering (sharedStrings order changed upstream; re-map displayed text) ---
$ws.Range("A88").Value = "Principado de Andorra"
$ws.Range("A89").Value = "Letonia"

$ws.Range("A164").Value = "San Martin (Parte Francesa)"
$ws.Range("A165").Value = "Benin"
$ws.Range("A166").Value = "Guam"
$ws.Range("A167").Value = "Sudan"

# --- Updated statistics ---
$ws.Range("B4").Value = 615406
$ws.Range("C4").Value = 1520
$ws.Range("D4").Value = 38879
$ws.Range("E4").Value = 550363
$ws.Range("G4").Value = 117
$ws.Range("H4").Value = 26164
$ws.Range("B8").Value = 132500
$ws.Range("C8").Value = 290
$ws.Range("E8").Value = 56379
$ws.Range("G8").Value = 26
$ws.Range("H8").Value = 3521
$ws.Range("B20").Value = 14325
$ws.Range("C20").Value = 99
$ws.Range("E20").Value = 5834
$ws.Range("F37").Value = 71
$ws.Range("F54").Value = 117
$ws.Range("B85").Value = 747
$ws.Range("C85").Value = 34
$ws.Range("E85").Value = 606
$ws.Range("B88").Value = 673
$ws.Range("C88").Value = 14
$ws.Range("D88").Value = 169
$ws.Range("E88").Value = 471
$ws.Range("F88").Value = 17
$ws.Range("G88").Value = 2
$ws.Range("H88").Value = 33
$ws.Range("B89").Value = 666
$ws.Range("C89").Value = 9
$ws.Range("D89").Value = 44
$ws.Range("E89").Value = 617
$ws.Range("F89").Value = 3
$ws.Range("H89").Value = 5
$ws.Range("D109").Value = 65
$ws.Range("E109").Value = 250
$ws.Range("C164").Value = 3
$ws.Range("D164").Value = 13
$ws.Range("E164").Value = 20
$ws.Range("F164").Value = 5
$ws.Range("H164").Value = 2
$ws.Range("B165").Value = 35
$ws.Range("D165").Value = 18
$ws.Range("E165").Value = 16
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 31
$ws.Range("H166").Value = 1
$ws.Range("D167").Value = 4
$ws.Range("E167").Value = 23
$ws.Range("F167").Value = 0
$ws.Range("H167").Value = 5

# --- Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 16:22"
